$wb = $excel.ActiveWorkbook

# --- Update the "Date" value on the Metadata sheet ---
$meta = $wb.Worksheets.Item("Metadata")
$meta.Range("B8").Value = "2025-10-02T18:31:12+01:00"

# --- Set "Case Sensitive" value (row 20, column B) to "true" ---
$meta.Range("B20").Value = "true"
